$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.561.85"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +5.20%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.725.18"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.98%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.34"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.32%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5390"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.65%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2703"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06619"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.86%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +5.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07759"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.79%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.658"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.13%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.744.87"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +5.31%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.961.97"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +3.94%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5896"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.62%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8319"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.97%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.14"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +3.68%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "27.570.02"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +5.35%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "225.60"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +16.94%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.765"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +1.80%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.27%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.132"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.45%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.005"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.02"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.40%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.699"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +11.66%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1238"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.94%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.435"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.61%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +5.02%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05604"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.26%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.39%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.604"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.40%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.474"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.31%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.671"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +6.23%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9657"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.99%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.825"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.50%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.448"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.79%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5954"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.82%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01655"

$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.898"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.59%  "

$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8593"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.77%  "

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.062.55"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.49%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.005"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.58"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.26%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.867.68"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.86%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈113"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +8.95%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "59.32"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.24%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.264"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.53%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4433"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.96%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.004"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.34%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05286"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.95%  "
